$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.840.70'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -5.04%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.211.76'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -6.33%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '315.93'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +1.54%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '99.40'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -8.09%  '
$ws.Range("E7").Value = '  -6.26%  '
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("E9").Value = '  -7.69%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.98'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -9.09%  '
$ws.Range("E11").Value = '  -2.77%  '
$ws.Range("E12").Value = '  -9.16%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.72'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -8.29%  '
$ws.Range("E14").Value = '  -2.50%  '
$ws.Range("B15").Value = 'Polygon'
$ws.Range("C15").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.859'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -11.54%  '
$ws.Range("B16").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C16").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.550.28'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -6.28%  '
$ws.Range("E17").Value = '  -6.56%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.214.30'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -6.64%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '42.759.97'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -5.21%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.63'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +2.61%  '
$ws.Range("E21").Value = '  -9.06%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.41'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -10.50%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '65.30'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -10.49%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.14'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -10.16%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '235.60'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -8.57%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.13'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -7.41%  '
$ws.Range("E27").Value = '  -0.04%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.00'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -9.40%  '
$ws.Range("E29").Value = '  -4.44%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.32'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -12.00%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0895'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -7.91%  '
$ws.Range("E32").Value = '  -7.90%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '34.34'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -7.16%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '155.47'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -7.38%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.79'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -7.88%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.15'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +9.13%  '
$ws.Range("E37").Value = '  +9.42%  '
$ws.Range("E38").Value = '  -6.56%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.90'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -1.53%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.42'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -4.99%  '
$ws.Range("E41").Value = '  -8.97%  '
$ws.Range("E42").Value = '  -7.52%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.904.74'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +0.60%  '
$ws.Range("E44").Value = '  -0.07%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '12.45'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -2.92%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '88.32'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -10.63%  '
$ws.Range("E47").Value = '  -8.97%  '
$ws.Range("B48").Value = 'MultiversX'
$ws.Range("C48").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '61.10'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -11.73%  '
$ws.Range("B49").Value = 'THORChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.36'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -4.63%  '
$ws.Range("B50").Value = 'ordi'
$ws.Range("C50").Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '76.67'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -8.12%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '102.68'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -6.35%  '
